$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 2, shifting existing data down
$ws.Rows.Item(2).Insert()

# The inserted row picked up stray formatting; clear it and
# re-apply the date style (from column A of the row below) to A2 only
$ws.Range("A2:E2").ClearFormats()
$ws.Range("A3").Copy()
$ws.Range("A2").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Write the full, corrected data set (all rows shifted + new values)
$ws.Range("A2").Value = 39400
$ws.Range("B2").Value = 2007
$ws.Range("C2").Value = 1.75539628881467
$ws.Range("D2").Value = 2008
$ws.Range("E2").Value = 1.327368416067398
$ws.Range("A3").Value = 39765
$ws.Range("B3").Value = 2008
$ws.Range("C3").Value = 2.213911448916162
$ws.Range("D3").Value = 2009
$ws.Range("E3").Value = 2.649257112350067
$ws.Range("A4").Value = 40130
$ws.Range("B4").Value = 2009
$ws.Range("C4").Value = 2.533533936850563
$ws.Range("D4").Value = 2010
$ws.Range("E4").Value = 1.815660192323709
$ws.Range("A5").Value = 40494
$ws.Range("B5").Value = 2010
$ws.Range("C5").Value = 2.088987486264915
$ws.Range("D5").Value = 2011
$ws.Range("E5").Value = 2.332261646026201
$ws.Range("A6").Value = 40862
$ws.Range("B6").Value = 2011
$ws.Range("C6").Value = 1.212544822741002
$ws.Range("D6").Value = 2012
$ws.Range("E6").Value = 1.839804681163293
$ws.Range("A7").Value = 41228
$ws.Range("B7").Value = 2012
$ws.Range("C7").Value = 1.196776590518644
$ws.Range("D7").Value = 2013
$ws.Range("E7").Value = 0.670590452940556
$ws.Range("A8").Value = 41592
$ws.Range("B8").Value = 2013
$ws.Range("C8").Value = 0.4712609263772594
$ws.Range("D8").Value = 2014
$ws.Range("E8").Value = 0.8520644823059476
$ws.Range("A9").Value = 41957
$ws.Range("B9").Value = 2014
$ws.Range("C9").Value = 0.8783377572271434
$ws.Range("D9").Value = 2015
$ws.Range("E9").Value = 1.474590898715178
$ws.Range("A10").Value = 42321
$ws.Range("B10").Value = 2015
$ws.Range("C10").Value = 2.29066283401107
$ws.Range("D10").Value = 2016
$ws.Range("E10").Value = 2.597902967862775
$ws.Range("A11").Value = 42689
$ws.Range("B11").Value = 2016
$ws.Range("C11").Value = 4.109890522944348
$ws.Range("D11").Value = 2017
$ws.Range("E11").Value = 3.628019428949036
$ws.Range("A12").Value = 43053
$ws.Range("B12").Value = 2017
$ws.Range("C12").Value = 1.336316831462692
$ws.Range("D12").Value = 2018
$ws.Range("E12").Value = 1.626630409005325
$ws.Range("A13").Value = 43418
$ws.Range("B13").Value = 2018
$ws.Range("C13").Value = 1.197912858979611
$ws.Range("D13").Value = 2019
$ws.Range("E13").Value = 1.216371234267344
$ws.Range("A14").Value = 43783
$ws.Range("B14").Value = 2019
$ws.Range("C14").Value = 1.727537197898665
$ws.Range("D14").Value = 2020
$ws.Range("E14").Value = 2.164378481800822
$ws.Range("A15").Value = 44159
$ws.Range("B15").Value = 2020
$ws.Range("C15").Value = 3.647228437274408
$ws.Range("D15").Value = 2021
$ws.Range("E15").Value = 3.845906281600109
$ws.Range("A16").Value = 44525
$ws.Range("B16").Value = 2021
$ws.Range("C16").Value = 2.777797690741424
$ws.Range("D16").Value = 2022
$ws.Range("E16").Value = 1.875884305456199
$ws.Range("A17").Value = 44890
$ws.Range("B17").Value = 2022
$ws.Range("C17").Value = 0.6994919452575576
$ws.Range("D17").Value = 2023
$ws.Range("E17").Value = -0.2388228654152447
$ws.Range("A18").Value = 45254
$ws.Range("B18").Value = 2023
$ws.Range("C18").Value = -1.432689847121871
$ws.Range("D18").Value = 2024
$ws.Range("E18").Value = -0.7896638887521124
$ws.Range("A19").Value = 45618
$ws.Range("B19").Value = 2024
$ws.Range("C19").Value = 2.033479419175133
$ws.Range("D19").Value = 2025
$ws.Range("E19").Value = 1.424898175306621

Write-Host "Done"
